# Update building block types in the Protein extraction template:
# - bump Version to 1.1.8
# - rename several "Parameter [...]" headers to "Component [...]" /
#   "Characteristic [...]" building blocks
# - fill in the previously-empty OBI:0000468 term source ref / accession
#   number for the "protein column" building block

$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump the Version value ---
$wsIsa = $wb.Worksheets.Item("isa_template")
$wsIsa.Range("B4").Value = "1.1.8"

# --- 2EXT02_Protein sheet: rename header / table column building blocks ---
$wsTable = $wb.Worksheets.Item("2EXT02_Protein")

$wsTable.Range("B1").Value = "Component [cleavage agent name]"
$wsTable.Range("E1").Value = "Characteristic [molecule]"
$wsTable.Range("H1").Value = "Characteristic [sample state]"
$wsTable.Range("N1").Value = "Component [extraction buffer]"
$wsTable.Range("X1").Value = "Component [protein column]"
$wsTable.Range("Y1").Value = "Term Source REF (OBI:0000468)"
$wsTable.Range("Z1").Value = "Term Accession Number (OBI:0000468)"
